# Helper: force a literal piece of text into a cell as a Text-typed value,
# bypassing Excel's automatic number/date/percentage inference (so "4564",
# "0", "1.06%" etc. are stored as text, not numbers).
function Set-TextValue {
    param($range, [string]$text)
    $range.Formula = '="' + $text + '"'
    $range.Copy()
    $range.PasteSpecial(-4163)
}

$wb = $excel.ActiveWorkbook

$wsBatting = $wb.Worksheets.Item("ODI Batting")
$wsBowling = $wb.Worksheets.Item("ODI Bowling")

# ---------------------------------------------------------------------------
# 1) "ODI Batting": the two placeholder/blank cells in column B (rows 2 & 3)
#    are removed entirely.
# ---------------------------------------------------------------------------
$wsBatting.Range("B2").ClearContents()
$wsBatting.Range("B3").ClearContents()

# ---------------------------------------------------------------------------
# 2) Add the new "ODI Batting Extra" sheet after "ODI Bowling" (so it becomes
#    the 4th sheet, sheetId 4).
# ---------------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($null, $wsBowling)
$newSheet.Name = "ODI Batting Extra"

# Bring over the bold / centered / bordered header style used by the other
# sheets' header rows (style index 1) by copying an existing header range
# and pasting only formats.
$wsBatting.Range("A1:F1").Copy()
$newSheet.Range("A1:F1").PasteSpecial(-4122)

# Header row
$newSheet.Range("A1").Value = "MATCH_CODE"
$newSheet.Range("B1").Value = "BATTING_POSITION"
$newSheet.Range("C1").Value = "NUM_4"
$newSheet.Range("D1").Value = "NUM_6"
$newSheet.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$newSheet.Range("F1").Value = "MAN_OF_MATCH"

# Row 2
Set-TextValue $newSheet.Range("A2") "4564"
$newSheet.Range("B2").Value = 11
Set-TextValue $newSheet.Range("F2") "NO"

# Row 3
Set-TextValue $newSheet.Range("A3") "4565"
Set-TextValue $newSheet.Range("F3") "NO"

# Row 4
Set-TextValue $newSheet.Range("A4") "4597"
$newSheet.Range("B4").Value = 9
Set-TextValue $newSheet.Range("C4") "0"
Set-TextValue $newSheet.Range("D4") "0"
Set-TextValue $newSheet.Range("E4") "1.06%"
Set-TextValue $newSheet.Range("F4") "NO"

# Restore the originally active sheet (the diff does not change bookViews).
$wb.Worksheets.Item(1).Activate()

Write-Output "edit applied"
